$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 103306942
$ws.Range("B2").Value = 78602
$ws.Range("E2").Value = 6463
$ws.Range("F2").Value = "Bårdlav"
$ws.Range("G2").Value = "Nephroma parile"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 615243.9621432471
$ws.Range("R2").Value = 7223727.877828724
$ws.Range("AJ2").Value = "sälg"
$ws.Range("AK2").Value = "Salix caprea"
$ws.Range("AO2").Value = "Salix caprea"
$ws.Range("A3").Value = 103306943
$ws.Range("B3").Value = 89388
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 1108
$ws.Range("F3").Value = "Harticka"
$ws.Range("G3").Value = "Pelloporus leporinus"
$ws.Range("H3").Value = "(Fr.) Krieglst."
$ws.Range("Q3").Value = 615238.5067762507
$ws.Range("R3").Value = 7223748.74792649
$ws.Range("AJ3").Value = "gran"
$ws.Range("AK3").Value = "Picea abies"
$ws.Range("AO3").Value = "Picea abies"
$ws.Range("A4").Value = 103306939
$ws.Range("B4").Value = 56395
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("Q4").Value = 615264.2335436697
$ws.Range("R4").Value = 7223629.999263954
$ws.Range("AJ4").Value = $null
$ws.Range("AK4").Value = $null
$ws.Range("AO4").Value = $null
$ws.Range("A5").Value = 103306940
$ws.Range("B5").Value = 96334
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("M5").Value = $null
$ws.Range("Q5").Value = 615261.0371939046
$ws.Range("R5").Value = 7223636.199837528
$ws.Range("A6").Value = 103306941
$ws.Range("B6").Value = 96237
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 220093
$ws.Range("F6").Value = "Korallrot"
$ws.Range("G6").Value = "Corallorhiza trifida"
$ws.Range("H6").Value = "Châtel."
$ws.Range("Q6").Value = 615248.9590524008
$ws.Range("R6").Value = 7223718.795936605
$ws.Range("A7").Value = 103306937
$ws.Range("B7").Value = 56411
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 100049
$ws.Range("F7").Value = "Spillkråka"
$ws.Range("G7").Value = "Dryocopus martius"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("M7").Value = "äldre spår"
$ws.Range("Q7").Value = 615266.3006582296
$ws.Range("R7").Value = 7223609.418203933
$ws.Range("A8").Value = 103306938
$ws.Range("B8").Value = 56395
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = "Tretåig hackspett"
$ws.Range("G8").Value = "Picoides tridactylus"
$ws.Range("Q8").Value = 615266.2679016123
$ws.Range("R8").Value = 7223610.260257882
$ws.Range("A9").Value = 103306948
$ws.Range("B9").Value = 90653
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 4364
$ws.Range("F9").Value = "Dropptaggsvamp"
$ws.Range("G9").Value = "Hydnellum ferrugineum"
$ws.Range("H9").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("M9").Value = $null
$ws.Range("Q9").Value = 615397.7041400182
$ws.Range("R9").Value = 7223725.850324323
